$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text (matching original inlineStr string cells)
# instead of auto-converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '38.811.36'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '2.103.71'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '227.25'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("D7").Value = '62.12'
$ws.Range("E7").Value = '  +3.01%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.95%  '
$ws.Range("D10").Value = '0.0844'
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '15.86'
$ws.Range("E12").Value = '  +5.90%  '
$ws.Range("D13").Value = '2.416.42'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '2.111.53'
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").Value = '38.805.04'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '71.78'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D21").Value = '0.0₃0845'
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").Value = '227.42'
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = '9.66'
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").Value = '170.61'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("D29").Value = '1.41'
$ws.Range("E29").Value = '  +2.44%  '
$ws.Range("D30").Value = '19.38'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("E31").Value = '  +8.92%  '
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.79'
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("B35").Value = 'THORChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D35").Value = '7.15'
$ws.Range("E35").Value = '  +13.55%  '
$ws.Range("D36").Value = '0.0613'
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D38").Value = '3.51'
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = '18.04'
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").Value = '0.0229'
$ws.Range("E41").Value = '  +3.49%  '
$ws.Range("D42").Value = '101.76'
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("D43").Value = '1.525.13'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("E44").Value = '  +8.04%  '
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.0913'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '7.75'
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +5.25%  '
$ws.Range("D49").Value = '4.16'
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("D51").Value = '2.304.15'
$ws.Range("E51").Value = '  +1.13%  '
